# Update existing rows 16-27: shift the "Fecha" (D) and volume/price values
# down one record (the newest two weekly records get inserted at the top,
# rows 16 and 19), and append two more historical rows (28 and 29) that
# previously sat at rows 26 and 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 ---
$ws.Range("D16").Value = 44846
$ws.Range("J16").Value = 7900

# --- Row 17 ---
$ws.Range("D17").Value = 44189
$ws.Range("J17").Value = 16000

# --- Row 18 ---
$ws.Range("D18").Value = 44160

# --- Row 19 ---
$ws.Range("D19").Value = 44845
$ws.Range("J19").Value = 7900

# --- Row 20 ---
$ws.Range("D20").Value = 44159
$ws.Range("J20").Value = 7000
$ws.Range("K20").Value = 3000
$ws.Range("M20").Value = 3000
$ws.Range("P20").Value = 30

# --- Row 21 ---
$ws.Range("D21").Value = 44215
$ws.Range("J21").Value = 16000

# --- Row 22 ---
$ws.Range("D22").Value = 44210
$ws.Range("J22").Value = 8800
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = 2750
$ws.Range("P22").Value = 28

# --- Row 23 ---
$ws.Range("D23").Value = 44166
$ws.Range("J23").Value = 7000

# --- Row 24 ---
$ws.Range("D24").Value = 44161

# --- Row 25 ---
$ws.Range("D25").Value = 44231
$ws.Range("J25").Value = 12000

# --- Row 26 ---
$ws.Range("D26").Value = 44162
$ws.Range("J26").Value = 7000

# --- Row 27 ---
$ws.Range("D27").Value = 44204
$ws.Range("J27").Value = 7000

# --- New row 28 (same record that used to be row 26) ---
$ws.Range("A28").Value = 6
$ws.Range("B28").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 44181
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 300000001
$ws.Range("G28").Value = "Rabanito"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 12000
$ws.Range("K28").Value = 3000
$ws.Range("L28").Value = 3000
$ws.Range("M28").Value = 3000
$ws.Range("N28").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O28").Value = "Provincia de Chacabuco"
$ws.Range("P28").Value = 30
$ws.Range("Q28").Value = 100
$ws.Range("R28").Value = "Hortaliza"

# --- New row 29 (same record that used to be row 27) ---
$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44187
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 300000001
$ws.Range("G29").Value = "Rabanito"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 12000
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = 3000
$ws.Range("N29").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O29").Value = "Provincia de Chacabuco"
$ws.Range("P29").Value = 30
$ws.Range("Q29").Value = 100
$ws.Range("R29").Value = "Hortaliza"
